# CVUmerHussain.docx edit script
# Applies (per the commit diff):
#  1. Phone number paragraph: "07584421701" -> "+44" / "7584421701",
#     bold (incl. complex-script bold) + complex-script size 12pt (24 half-pts).
#  2. "Finance: C<tab><tab>Business Studies: MP" -> split off the leading
#     "B" of "Business" into its own run (re-typed in place).
#  3. Three similar "<tab><spaces>" runs (end of GCSE lines) get the
#     trailing spaces split into their own run.
#  4. bookmarkStart around "_GoBack" is left semantically the same.
#  5. Page setup: explicit Portrait orientation on the section.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Phone number: 07584421701 -> +44 7584421701 (typed as two runs)
# ---------------------------------------------------------------------
$phone = $d.Paragraphs(4).Range
$phone.Font.Bold = $true
$phone.Font.BoldBi = $true
$phone.Font.SizeBi = 12

$found = $d.Content.Find.Execute("07584421701", $false, $false, $false, $false, $false, $true, 1, $false, "+447584421701", 2)

# Re-touch the "+44" prefix so it keeps living in its own run (mirrors the
# author retyping the country code in front of the existing number).
$prefix = $d.Paragraphs(4).Range
$prefixStart = $prefix.Start
$plus = $d.Range($prefixStart, $prefixStart + 3)
$plus.Font.Size = 99
$plus.Font.Size = 12

# ---------------------------------------------------------------------
# 2) "...Finance: C<tab><tab>Business Studies: MP" -> split the "B"
#    of "Business" into its own run.
# ---------------------------------------------------------------------
$rngB = $d.Content
$null = $rngB.Find.Execute("usiness Studies", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$bRange = $d.Range($rngB.Start - 1, $rngB.Start)
$bRange.Font.Size = 99
$bRange.Font.Size = 12

# ---------------------------------------------------------------------
# 3) Split the trailing-space runs off their preceding <w:tab/> on the
#    three GCSE lines (English Literature / English Language / Geography).
# ---------------------------------------------------------------------
function Split-TrailingSpaces($afterText, $spaceCount) {
    $r = $d.Content
    $ok = $r.Find.Execute($afterText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($ok) {
        $spacesStart = $r.End + 1
        $spaces = $d.Range($spacesStart, $spacesStart + $spaceCount)
        $spaces.Font.Size = 99
        $spaces.Font.Size = 12
    }
}

Split-TrailingSpaces "English Literature: B" 3
Split-TrailingSpaces "English Language: C" 2
Split-TrailingSpaces "Geography: C" 6

# ---------------------------------------------------------------------
# 4) GCSE / _GoBack bookmark region - no textual change, left as-is.
# ---------------------------------------------------------------------

# ---------------------------------------------------------------------
# 5) Section page setup: explicit portrait orientation.
# ---------------------------------------------------------------------
$d.PageSetup.Orientation = 0
